# "update radar chart and RAD relative import"
#
# The underlying edit is a re-vote on the "RAD" (Votare sulla base di: IPT)
# column for the second evaluator block - several "Voti IPT" scores for the
# F column (rows 39,44,48,53,55,60,64) were lowered. These feed the
# AVERAGE() formulas in F71/F72/F73/F75 and, in turn, the radar chart
# (chart2.xml) that plots column F (F67:F73) per category - hence "update
# radar chart" in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core data edits: lower several "Voti IPT" (col F) scores ---
$ws.Range("F39").Value = 2   # was 3
$ws.Range("F44").Value = 1   # was 2
$ws.Range("F48").Value = 2   # was 4
$ws.Range("F53").Value = 2   # was 3
$ws.Range("F55").Value = 1   # was 2
$ws.Range("F60").Value = 0   # was 4
$ws.Range("F64").Value = 1   # was 4

# --- Recalculate so the AVERAGE() rollups (F71/F72/F73/F75) and the
#     dependent radar chart series pick up the new figures ---
$excel.CalculateFullRebuild()

# --- Refresh the radar charts so their cached plot data is rebuilt from
#     the new worksheet values ---
$ws.ChartObjects().Item(1).Chart.Refresh() | Out-Null
$ws.ChartObjects().Item(2).Chart.Refresh() | Out-Null

# --- Restore the view/selection the author ended up with after scrolling
#     down to review rows 37+ and selecting the newly-recomputed F69:I69
#     rollup cell ---
$ws.Range("A37").Select() | Out-Null
$ws.Range("F69:I69").Select() | Out-Null
